# Generate Report for Handoff
# Update status & timestamps across Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-29-19 08:29:45"

# zh-cn sheet: Status column C, Latest Handoff Datetime column E
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-19 08:29:42"

# de-de sheet: Status column C, Latest Handoff Datetime column E
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-19 08:29:45"
